$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 439, holds the "Förändrad" (changed) date.
# All of these cells currently hold serial date 45186 (2023-09-17) and
# must be updated to 45188 (2023-09-19).
$ws.Range("C2:C439").Value = 45188
